# Correct for res/com sf bug
# Updates the "res" (I/J columns) and "com" (L/M columns) figures on the
# SummaryTable_County_wTotals sheet for each county row. UnitsDensity pulls
# these numbers in via formulas, so it recalculates automatically.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("SummaryTable_County_wTotals")

$src.Range("I2").Value = 84100
$src.Range("J2").Value = 74000
$src.Range("L2").Value = 57
$src.Range("M2").Value = 47

$src.Range("I3").Value = 109200
$src.Range("J3").Value = 98900
$src.Range("L3").Value = 77
$src.Range("M3").Value = 68

$src.Range("I4").Value = 365800
$src.Range("J4").Value = 320800
$src.Range("L4").Value = 86
$src.Range("M4").Value = 67

$src.Range("I5").Value = 144800
$src.Range("J5").Value = 126200
$src.Range("L5").Value = 84
$src.Range("M5").Value = 68

$src.Range("I6").Value = 703900
$src.Range("J6").Value = 619900
$src.Range("L6").Value = 80
$src.Range("M6").Value = 64
